# repull data, push all data, mean calculation
# Updates the "dSF" (F) column with freshly re-pulled delta values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -2
    4  = -2
    5  = 2
    7  = 1
    8  = -2
    9  = -1
    10 = -4
    11 = -1
    12 = 7
    13 = 0
    16 = -1
    17 = 1
    19 = -3
    20 = -3
    21 = -1
    22 = -6
    24 = 4
    25 = -4
    26 = 7
    27 = 1
    28 = -3
    29 = 6
    30 = 2
    31 = -4
    32 = 1
    33 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
